$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Header: BTec_Logo-Orange picture, rename image1.jpg -> image2.jpg ---
for ($hi = 1; $hi -le 3; $hi++) {
    $hdr = $sec.Headers.Item($hi)
    if ($hdr.Exists) {
        $cnt = $hdr.Range.InlineShapes.Count
        for ($i = 1; $i -le $cnt; $i++) {
            $shp = $hdr.Range.InlineShapes.Item($i)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image2.jpg"
            }
        }
    }
}

# --- Footers: Pearson logo pictures, rename image2.png -> image1.png ---
for ($fi = 1; $fi -le 3; $fi++) {
    $ftr = $sec.Footers.Item($fi)
    if ($ftr.Exists) {
        $cnt = $ftr.Range.InlineShapes.Count
        for ($i = 1; $i -le $cnt; $i++) {
            $shp = $ftr.Range.InlineShapes.Item($i)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shp.Name = "image1.png"
            }
        }
    }
}
